# Update Work Week and Social Spending
#
# The source workbook ("Data" sheet) holds one GDP-per-Capita observation
# per row for Oman (Country Code 512), one row per year, 1950 onward.
# Column layout: A=Country Code, B=Country Name, C=Indicator, D=Year,
# E=Value (stored as text, like the rest of this workbook).
#
# This refresh replaces every existing yearly value (1950-2010, rows 2-62)
# with the newer published figures, and appends six new years (2011-2016,
# rows 63-68) that didn't exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated GDP per Capita values for years 1950..2016, in order.
$values = @(
    "993","1036","1079","1127","1173","1221","1274","1325","1380","1435",
    "1490","1471","1728","1758","1714","1678","1722","2788","4932","6028",
    "6056","5923","6271","5227","5644","6802","7328","6998","6513","6446",
    "6491","7210","7651","8521","9526","10433","10268","10699","10633","10692",
    "10327","11131.0001614999","12251.1576128381","13226.6432323667","14079.5330556326","15292.1759621806","16491.1678287698","18518.4805057695","20213.0203779339","21519.6708744851",
    "24186.3653983285","26431.1005102425","27141.8429085551","27313.979688054","28607.0843640724","30392.5600695084","33374.8938759397","36424.9793543135","40909.9379307961","44244.6465509869",
    "46182.1601222189","44454","44008","41680","39418","38738",
    "38515"
)

$firstYear = 1950
$lastExistingRow = 62

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $year = $firstYear + $i
    $value = $values[$i]

    if ($row -gt $lastExistingRow) {
        # These rows don't exist yet (years 2011-2016) - create them in full.
        $ws.Cells.Item($row, 1).Value = 512
        $ws.Cells.Item($row, 2).Value = "Oman"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    # Force text storage (matches the rest of column E in this sheet) -
    # otherwise Excel auto-detects these numeric-looking strings as numbers.
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 5).ClearFormats()
}

"done"
